$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-13 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-14 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("22×58=1276", $true, $false, $false, $false, $false, $true, 1, $false, "34×83=2822", 2) | Out-Null
$d.Content.Find.Execute("64×95=6080", $true, $false, $false, $false, $false, $true, 1, $false, "86×99=8514", 2) | Out-Null
$d.Content.Find.Execute("56×27=1512", $true, $false, $false, $false, $false, $true, 1, $false, "68×73=4964", 2) | Out-Null
$d.Content.Find.Execute("57×59=3363", $true, $false, $false, $false, $false, $true, 1, $false, "78×25=1950", 2) | Out-Null
$d.Content.Find.Execute("84×81=6804", $true, $false, $false, $false, $false, $true, 1, $false, "87×52=4524", 2) | Out-Null
$d.Content.Find.Execute("30×66=1980", $true, $false, $false, $false, $false, $true, 1, $false, "78×70=5460", 2) | Out-Null
$d.Content.Find.Execute("41×40=1640", $true, $false, $false, $false, $false, $true, 1, $false, "37×37=1369", 2) | Out-Null
$d.Content.Find.Execute("80×93=7440", $true, $false, $false, $false, $false, $true, 1, $false, "35×16=560", 2) | Out-Null
$d.Content.Find.Execute("71×81=5751", $true, $false, $false, $false, $false, $true, 1, $false, "59×29=1711", 2) | Out-Null
$d.Content.Find.Execute("40×69=2760", $true, $false, $false, $false, $false, $true, 1, $false, "49×89=4361", 2) | Out-Null
$d.Content.Find.Execute("99×82=8118", $true, $false, $false, $false, $false, $true, 1, $false, "47×84=3948", 2) | Out-Null
$d.Content.Find.Execute("31×62=1922", $true, $false, $false, $false, $false, $true, 1, $false, "48×98=4704", 2) | Out-Null
$d.Content.Find.Execute("29×45=1305", $true, $false, $false, $false, $false, $true, 1, $false, "14×21=294", 2) | Out-Null
$d.Content.Find.Execute("58×43=2494", $true, $false, $false, $false, $false, $true, 1, $false, "22×16=352", 2) | Out-Null
$d.Content.Find.Execute("43×33=1419", $true, $false, $false, $false, $false, $true, 1, $false, "30×25=750", 2) | Out-Null
$d.Content.Find.Execute("48×32=1536", $true, $false, $false, $false, $false, $true, 1, $false, "38×77=2926", 2) | Out-Null
$d.Content.Find.Execute("36×63=2268", $true, $false, $false, $false, $false, $true, 1, $false, "47×69=3243", 2) | Out-Null
$d.Content.Find.Execute("78×48=3744", $true, $false, $false, $false, $false, $true, 1, $false, "75×18=1350", 2) | Out-Null
$d.Content.Find.Execute("30×49=1470", $true, $false, $false, $false, $false, $true, 1, $false, "43×45=1935", 2) | Out-Null
$d.Content.Find.Execute("71×15=1065", $true, $false, $false, $false, $false, $true, 1, $false, "54×67=3618", 2) | Out-Null
$d.Content.Find.Execute("68×11=748", $true, $false, $false, $false, $false, $true, 1, $false, "97×27=2619", 2) | Out-Null
$d.Content.Find.Execute("15×90=1350", $true, $false, $false, $false, $false, $true, 1, $false, "78×61=4758", 2) | Out-Null
$d.Content.Find.Execute("13×87=1131", $true, $false, $false, $false, $false, $true, 1, $false, "12×17=204", 2) | Out-Null
$d.Content.Find.Execute("67×70=4690", $true, $false, $false, $false, $false, $true, 1, $false, "69×69=4761", 2) | Out-Null
$d.Content.Find.Execute("35×48=1680", $true, $false, $false, $false, $false, $true, 1, $false, "90×36=3240", 2) | Out-Null
